$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new work-log entry for row 8
$ws.Range("A8").Value2 = 43452
$ws.Range("B8").Value2 = 'Continuation on work for "Backgrounds and Methods"'
$ws.Range("C8").Value2 = 1

# Copy the date formatting from the row above (A2) so the new date cell
# reuses the existing date style instead of creating a new one
$ws.Range("A2").Copy()
$ws.Range("A8").PasteSpecial(-4122)

# Move the active selection, matching the saved workbook view state
$ws.Range("M9").Select()
